$d = $word.ActiveDocument

# The diff removes three whole paragraphs that immediately follow the
# "LOQ4038: Química Orgânica II (Requisito fraco)" paragraph:
#   1. an empty paragraph
#   2. "Ver no Jupiter Salvar em pdf Salvar em docx"
#   3. "© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github
#       pages. Original theme under Creative Commons Attribution"
# The "LOQ4038..." paragraph itself, and the empty paragraph that comes
# right after those three (just before the page-break paragraph), are
# kept untouched.

# Locate the "LOQ4038..." paragraph using Find so we don't rely on
# brittle, hard-coded paragraph numbers.
$anchor = $d.Content
$null = $anchor.Find.Execute("LOQ4038: Química Orgânica II (Requisito fraco)", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Start -le $anchor.Start -and $p.Range.End -ge $anchor.End) {
        $anchorIndex = $i
        break
    }
}

if ($anchorIndex -eq -1) {
    throw "Could not locate the LOQ4038 paragraph"
}

# Delete the three paragraphs that directly follow the anchor paragraph.
$startDelete = $d.Paragraphs.Item($anchorIndex + 1).Range.Start
$endDelete = $d.Paragraphs.Item($anchorIndex + 3).Range.End

$delRange = $d.Range($startDelete, $endDelete)
$delRange.Delete()
